$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new column D (shifts old D "Duración..." data to E) ---
$ws.Columns.Item(4).Insert()

# --- Header row (row 1) text ---
$ws.Range("B1").Value = "Promedio de días lluviosos"
$ws.Range("C1").Value = "Cantidad promedio de rachas"
$ws.Range("D1").Value = "Cantidad promedio de rachas de más de un día"
$ws.Range("E1").Value = "Duración media de las rachas (Todas)"
$ws.Range("F1").Value = "Probabilidad media de días lluviosos"
$ws.Range("G1").Value = "Probabilidad media de días lluviosos si llovió el anterior"
$ws.Range("H1").Value = "Coeficiente de Besson"
$ws.Range("I1").Value = "Índice de irregularidad temporal"

# --- New column D values (Cantidad promedio de rachas de más de un día) ---
$ws.Range("D2").Value = 0.7837078651685393
$ws.Range("D3").Value = 0.6507042253521127
$ws.Range("D4").Value = 0.62395543175487467
$ws.Range("D5").Value = 0.74864864864864866
$ws.Range("D6").Value = 1.9146341463414629
$ws.Range("D7").Value = 2.167865707434053
$ws.Range("D8").Value = 1.554479418886199
$ws.Range("D9").Value = 2.1490384615384621
$ws.Range("D10").Value = 2.504807692307693
$ws.Range("D11").Value = 1.9928057553956831
$ws.Range("D12").Value = 1.2424999999999999
$ws.Range("D13").Value = 0.81481481481481477

# --- New columns F-I values ---
$ws.Range("F2").Value = 0.1249377038782169
$ws.Range("G2").Value = 0.21672209656760219
$ws.Range("H2").Value = 0.1114194470611314
$ws.Range("I2").Value = 0.26219601584003538
$ws.Range("F3").Value = 0.1159268253197346
$ws.Range("G3").Value = 0.2030538788285268
$ws.Range("H3").Value = 0.1040258428453498
$ws.Range("I3").Value = 0.26080985673628843
$ws.Range("F4").Value = 0.11545848683619379
$ws.Range("G4").Value = 0.1609563602599815
$ws.Range("H4").Value = 0.05825861471461286
$ws.Range("I4").Value = 0.27998156310202882
$ws.Range("F5").Value = 0.12926765475152599
$ws.Range("G5").Value = 0.2222877422877424
$ws.Range("H5").Value = 0.11436898014038389
$ws.Range("I5").Value = 0.33995218365974028
$ws.Range("F6").Value = 0.26866148701809561
$ws.Range("G6").Value = 0.46059140323544118
$ws.Range("H6").Value = 0.27413728988992159
$ws.Range("I6").Value = 0.76420218783685578
$ws.Range("F7").Value = 0.32451716046517642
$ws.Range("G7").Value = 0.41498451421228832
$ws.Range("H7").Value = 0.14224061283339651
$ws.Range("I7").Value = 0.94506850633060169
$ws.Range("F8").Value = 0.2374297039756309
$ws.Range("G8").Value = 0.27954130429383262
$ws.Range("H8").Value = 0.06073543911946271
$ws.Range("I8").Value = 0.6320988037655606
$ws.Range("F9").Value = 0.30032196675971912
$ws.Range("G9").Value = 0.37149677060993369
$ws.Range("H9").Value = 0.1089652582711768
$ws.Range("I9").Value = 0.7859253995556954
$ws.Range("F10").Value = 0.34639565685519541
$ws.Range("G10").Value = 0.41132927609497538
$ws.Range("H10").Value = 0.1037483005218994
$ws.Range("I10").Value = 0.9416550972530815
$ws.Range("F11").Value = 0.3012111278718963
$ws.Range("G11").Value = 0.39657771763109212
$ws.Range("H11").Value = 0.14854759191631159
$ws.Range("I11").Value = 0.80626652123412956
$ws.Range("F12").Value = 0.2046639784946237
$ws.Range("G12").Value = 0.3381940749846129
$ws.Range("H12").Value = 0.18109996176665821
$ws.Range("I12").Value = 0.48407981080511819
$ws.Range("F13").Value = 0.13555739362190969
$ws.Range("G13").Value = 0.2044166645109991
$ws.Range("H13").Value = 0.08981287210589134
$ws.Range("I13").Value = 0.28141457297832201

# --- Copy the existing bold/border header style onto the new header cells ---
$ws.Range("B1").Copy()
$ws.Range("F1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Wrap header text and set the taller row height ---
$ws.Range("B1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 75

# --- Number formats for the data columns ---
$ws.Range("B2:E13").NumberFormat = "0.00"
$ws.Range("H2:I13").NumberFormat = "0.00"
$ws.Range("F2:G13").NumberFormat = "0.00%"

# --- Column widths ---
$ws.Range("A1:I13").ColumnWidth = 14.7109375
